# Atualização de bases das ligas, do dia: 08-04-2024 às 21:28
#
# 1) Swap the betting-odds data (columns B:AC) between three pairs of rows
#    (the "id" column A stays put per row number).
# 2) Append two new match rows (150, 151) at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues($row, $values) {
    $cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "$($cols[$i])$row"
        $ws.Range($addr).Value2 = $values[$i]
    }
}

function Set-NewRowValues($row, $values) {
    $cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "$($cols[$i])$row"
        $ws.Range($addr).Value2 = $values[$i]
    }
}

# --- Swap rows 9 <-> 10 ---
$row9vals  = @(6865285,"Bosnia Herzegovina Premier Liga","Bosnia  Herzegovina Premier Liga",45150.5,"NK Igman Konjic","Sloga",1,0,"H",2,3.4,3.2,1.909,3.5,3.4,-0.5,1.95,1.85,2.5,1.85,1.95,0.909,-1,-1,0.95,-1,-1,0.95)
$row10vals = @(6865281,"Bosnia Herzegovina Premier Liga","Bosnia  Herzegovina Premier Liga",45150.5,"GOSK Gabela","Zvijezda 09",2,0,"H",1.75,4,3.5,1.75,4,3.4,-0.5,1.8,2,2.5,1.85,1.95,0.75,-1,-1,0.8,-1,-1,0.95)
Set-RowValues 9 $row9vals
Set-RowValues 10 $row10vals

# --- Swap rows 36 <-> 37 ---
$row36vals = @(6865299,"Bosnia Herzegovina Premier Liga","Bosnia  Herzegovina Premier Liga",45186.61458333334,"Siroki Brijeg","Zvijezda 09",2,1,"H",1.25,5.5,8,1.4,4.75,5.75,-1.25,1.9,1.9,2.75,1.85,1.95,0.3999999999999999,-1,-1,-0.5,0.45,0.425,-0.5)
$row37vals = @(6864629,"Bosnia Herzegovina Premier Liga","Bosnia  Herzegovina Premier Liga",45186.61458333334,"Borac Banja Luka","NK Posusje",1,0,"H",1.363,4.5,6.5,1.363,4.2,6.5,-1.25,1.95,1.85,2.5,1.925,1.875,0.363,-1,-1,-0.5,0.425,-1,0.875)
Set-RowValues 36 $row36vals
Set-RowValues 37 $row37vals

# --- Swap rows 49 <-> 50 ---
$row49vals = @(6865311,"Bosnia Herzegovina Premier Liga","Bosnia  Herzegovina Premier Liga",45200.41666666666,"Sloga","GOSK Gabela",3,2,"H",1.833,3.6,3.4,1.909,3.4,3.3,-0.5,1.925,1.875,2.25,1.825,1.975,0.909,-1,-1,0.925,-1,0.825,-1)
$row50vals = @(6865310,"Bosnia Herzegovina Premier Liga","Bosnia  Herzegovina Premier Liga",45200.41666666666,"NK Igman Konjic","Zrinjski Mostar",0,2,"A",3.4,3.6,1.833,4.75,4.75,1.45,1.25,1.775,2.025,2.75,1.85,1.95,-1,-1,0.45,-1,1.025,-1,0.95)
Set-RowValues 49 $row49vals
Set-RowValues 50 $row50vals

# --- Append new rows 150 and 151 ---
# Copy formatting from the last existing data row (149) for the id column (A, bold/border)
# and the date column (E, custom date/time format), then fill in values.
$ws.Range("A149").Copy()
$ws.Range("A150:A151").PasteSpecial(-4122)
$ws.Range("E149").Copy()
$ws.Range("E150:E151").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$row150vals = @(148,7952744,"Bosnia Herzegovina Premier Liga","Bosnia  Herzegovina Premier Liga",45390.4375,"Siroki Brijeg","GOSK Gabela",2,2,"D",1.727,3.25,4.333,1.65,3.5,4.5,-0.75,1.875,1.925,2.5,1.925,1.875,-1,2.5,-1,-1,0.925,0.925,-1)
$row151vals = @(149,7952457,"Bosnia Herzegovina Premier Liga","Bosnia  Herzegovina Premier Liga",45390.54166666666,"Zrinjski Mostar","Zvijezda 09",4,0,"H",1.166,6.5,10,1.166,6,15,-2,1.9,1.9,3,1.875,1.925,0.1659999999999999,-1,-1,0.8999999999999999,-1,0.875,-1)
Set-NewRowValues 150 $row150vals
Set-NewRowValues 151 $row151vals

Write-Output "done"
